# Natmi following Dr Hou advice
# Update the LR-pair stats for rows 2-4 (new detection counts / expression
# values) and add a new "M2" target-cluster row, re-using the "sCs" label
# that row 4 used to carry (row 4 becomes the "M2" cluster instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (Target cluster: ECs) ----
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.621828666666667
$ws.Range("H2").Value = 4.865486
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.451301666666667
$ws.Range("N2").Value = 7.353905
$ws.Range("O2").Value = 0.1099648918267336
$ws.Range("P2").Value = 0.1099648918267337
$ws.Range("Q2").Value = 3.975591313647778
$ws.Range("R2").Value = 35.78032182283
$ws.Range("S2").Value = 0.1099648918267336
$ws.Range("T2").Value = 0.1099648918267337

# ---- Row 3 (Target cluster: FAPs) ----
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.621828666666667
$ws.Range("H3").Value = 4.865486
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.637306
$ws.Range("N3").Value = 43.911918
$ws.Range("O3").Value = 0.6566265559283671
$ws.Range("P3").Value = 0.6566265559283672
$ws.Range("Q3").Value = 23.739202473572
$ws.Range("R3").Value = 213.652822262148
$ws.Range("S3").Value = 0.6566265559283671
$ws.Range("T3").Value = 0.6566265559283672

# ---- Row 4 (Target cluster changes from sCs -> M2) ----
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.621828666666667
$ws.Range("H4").Value = 4.865486
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02548533333333333
$ws.Range("N4").Value = 0.076456
$ws.Range("O4").Value = 0.001143266845234572
$ws.Range("P4").Value = 0.001143266845234572
$ws.Range("Q4").Value = 0.04133284417955555
$ws.Range("R4").Value = 0.371995597616
$ws.Range("S4").Value = 0.001143266845234572
$ws.Range("T4").Value = 0.001143266845234572

# ---- Row 5 (NEW row, Target cluster: sCs) ----
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.621828666666667
$ws.Range("H5").Value = 4.865486
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.177582333333334
$ws.Range("N5").Value = 15.532747
$ws.Range("O5").Value = 0.2322652853996647
$ws.Range("P5").Value = 0.2322652853996648
$ws.Range("Q5").Value = 8.397151452226888
$ws.Range("R5").Value = 75.574363070042
$ws.Range("S5").Value = 0.2322652853996647
$ws.Range("T5").Value = 0.2322652853996648
